$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = -21.947
$ws.Range("A18").Value = -22.38170000000002
$ws.Range("A20").Value = -20.63879999999998
$ws.Range("A27").Value = -21.78339999999998
$ws.Range("A69").Value = -21.63949999999999
$ws.Range("A76").Value = -19.87139999999998
$ws.Range("A82").Value = -21.9675
